$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,10
$row2[0,0] = 9.960511441256548
$row2[0,1] = -7.953374226292953
$row2[0,2] = 0.2335353073468205
$row2[0,3] = 0.8732301461903136
$row2[0,4] = -0.9276358786007236
$row2[0,5] = -1.886570575107837
$row2[0,6] = 0.9471476497731522
$row2[0,7] = -0.7219570139330117
$row2[0,8] = 0.298215423810231
$row2[0,9] = -1.41557403308908
$ws.Range("B2:K2").Value = $row2

$row3 = New-Object 'object[,]' 1,10
$row3[0,0] = -10.20317706590435
$row3[0,1] = -2.016267532264572
$row3[0,2] = -1.376572693421079
$row3[0,3] = -3.177438718212116
$row3[0,4] = -4.136373414719229
$row3[0,5] = -1.30265518983824
$row3[0,6] = -2.971759853544404
$row3[0,7] = -1.951587415801161
$row3[0,8] = -3.665376872700473
$row3[0,9] = -0.8370214991015008
$ws.Range("B3:K3").Value = $row3

$row4 = New-Object 'object[,]' 1,10
$row4[0,0] = -9.485417863121866
$row4[0,1] = -8.84572302427837
$row4[0,2] = -10.64658904906941
$row4[0,3] = -11.60552374557652
$row4[0,4] = -8.771805520695533
$row4[0,5] = -10.4409101844017
$row4[0,6] = -9.420737746658455
$row4[0,7] = -11.13452720355777
$row4[0,8] = -8.306171829958794
$row4[0,9] = -10.54150442982491
$ws.Range("B4:K4").Value = $row4

$row5 = New-Object 'object[,]' 1,10
$row5[0,0] = 9.325538775859098
$row5[0,1] = 7.52467275106806
$row5[0,2] = 6.565738054560947
$row5[0,3] = 9.399456279441935
$row5[0,4] = 7.730351615735771
$row5[0,5] = 8.750524053479014
$row5[0,6] = 7.036734596579703
$row5[0,7] = 9.865089970178675
$row5[0,8] = 7.629757370312555
$row5[0,9] = 8.981092428230218
$ws.Range("B5:K5").Value = $row5

$row6 = New-Object 'object[,]' 1,10
$row6[0,0] = -0.6149987959110895
$row6[0,1] = -1.573933492418202
$row6[0,2] = 1.259784732462786
$row6[0,3] = -0.4093199312433776
$row6[0,4] = 0.6108525064998651
$row6[0,5] = -1.102936950399446
$row6[0,6] = 1.725418423199526
$row6[0,7] = -0.5099141766665937
$row6[0,8] = 0.8414208812510687
$row6[0,9] = 0.2199829514341669
$ws.Range("B6:K6").Value = $row6

$row7 = New-Object 'object[,]' 1,10
$row7[0,0] = -3.721503964217595
$row7[0,1] = -0.8877857393366061
$row7[0,2] = -2.55689040304277
$row7[0,3] = -1.536717965299527
$row7[0,4] = -3.250507422198839
$row7[0,5] = -0.4221520485998669
$row7[0,6] = -2.657484648465986
$row7[0,7] = -1.306149590548324
$row7[0,8] = -1.927587520365226
$row7[0,9] = -1.427026823174395
$ws.Range("B7:K7").Value = $row7

$row8 = New-Object 'object[,]' 1,10
$row8[0,0] = 2.552244247186724
$row8[0,1] = 0.8831395834805599
$row8[0,2] = 1.903312021223803
$row8[0,3] = 0.1895225643244911
$row8[0,4] = 3.017877937923463
$row8[0,5] = 0.7825453380573438
$row8[0,6] = 2.133880395975006
$row8[0,7] = 1.512442466158104
$row8[0,8] = 2.013003163348936
$row8[0,9] = 1.345095091002794
$ws.Range("B8:K8").Value = $row8

$row9 = New-Object 'object[,]' 1,10
$row9[0,0] = 0.9443137618947609
$row9[0,1] = 1.964486199638004
$row9[0,2] = 0.2506967427386921
$row9[0,3] = 3.079052116337664
$row9[0,4] = 0.8437195164715449
$row9[0,5] = 2.195054574389207
$row9[0,6] = 1.573616644572305
$row9[0,7] = 2.074177341763137
$row9[0,8] = 1.406269269416995
$row9[0,9] = 1.835270244654998
$ws.Range("B9:K9").Value = $row9

$row10 = New-Object 'object[,]' 1,10
$row10[0,0] = -0.07696973751825081
$row10[0,1] = -1.790759194417562
$row10[0,2] = 1.03759617918141
$row10[0,3] = -1.19773642068471
$row10[0,4] = 0.1535986372329528
$row10[0,5] = -0.467839292583949
$row10[0,6] = 0.032721404606882
$row10[0,7] = -0.6351866677392595
$row10[0,8] = -0.2061856925012563
$row10[0,9] = -0.1485141439230462
$ws.Range("B10:K10").Value = $row10

$row11 = New-Object 'object[,]' 1,10
$row11[0,0] = -1.208164657383921
$row11[0,1] = 1.620190716215051
$row11[0,2] = -0.6151418836510686
$row11[0,3] = 0.7361931742665938
$row11[0,4] = 0.114755244449692
$row11[0,5] = 0.615315941640523
$row11[0,6] = -0.05259213070561841
$row11[0,7] = 0.3764088445323847
$row11[0,8] = 0.4340803931105948
$row11[0,9] = 0.5580467639488803
$ws.Range("B11:K11").Value = $row11

$row12 = New-Object 'object[,]' 1,10
$row12[0,0] = 1.524515675405693
$row12[0,1] = -0.7108169244604263
$row12[0,2] = 0.6405181334572361
$row12[0,3] = 0.01908020364033419
$row12[0,4] = 0.5196409008311652
$row12[0,5] = -0.1482671715149762
$row12[0,6] = 0.2807338037230269
$row12[0,7] = 0.338405352301237
$row12[0,8] = 0.4623717231395225
$row12[0,9] = -0.4310464000952693
$ws.Range("B12:K12").Value = $row12

$row13 = New-Object 'object[,]' 1,10
$row13[0,0] = -0.4503978874617036
$row13[0,1] = 0.9009371704559588
$row13[0,2] = 0.279499240639057
$row13[0,3] = 0.780059937829888
$row13[0,4] = 0.1121518654837466
$row13[0,5] = 0.5411528407217497
$row13[0,6] = 0.5988243892999598
$row13[0,7] = 0.7227907601382453
$row13[0,8] = -0.1706273630965465
$row13[0,9] = 0.5995033638472159
$ws.Range("B13:K13").Value = $row13

$row14 = New-Object 'object[,]' 1,10
$row14[0,0] = 0.06548217429746761
$row14[0,1] = -0.5559557555194342
$row14[0,2] = -0.0553950583286032
$row14[0,3] = -0.7233031306747446
$row14[0,4] = -0.2943021554367415
$row14[0,5] = -0.2366306068585314
$row14[0,6] = -0.1126642360202459
$row14[0,7] = -1.006082359255038
$row14[0,8] = -0.2359516323112753
$row14[0,9] = -0.5018523531907899
$ws.Range("B14:K14").Value = $row14

$row15 = New-Object 'object[,]' 1,9
$row15[0,0] = -0.4558763956168127
$row15[0,1] = 0.04468430157401831
$row15[0,2] = -0.6232237707721231
$row15[0,3] = -0.19422279553412
$row15[0,4] = -0.1365512469559099
$row15[0,5] = -0.01258487611762438
$row15[0,6] = -0.9060029993524162
$row15[0,7] = -0.1358722724086538
$row15[0,8] = -0.4017729932881683
$ws.Range("B15:J15").Value = $row15
$ws.Range("K15").ClearContents()

$row16 = New-Object 'object[,]' 1,8
$row16[0,0] = 0.3984708076537146
$row16[0,1] = -0.2694372646924268
$row16[0,2] = 0.1595637105455762
$row16[0,3] = 0.2172352591237863
$row16[0,4] = 0.3412016299620719
$row16[0,5] = -0.55221649327272
$row16[0,6] = 0.2179142336710425
$row16[0,7] = -0.04798648720847212
$ws.Range("B16:I16").Value = $row16
$ws.Range("J16").ClearContents()

$row17 = New-Object 'object[,]' 1,7
$row17[0,0] = -0.4176146551248722
$row17[0,1] = 0.0113863201131309
$row17[0,2] = 0.069057868691341
$row17[0,3] = 0.1930242395296265
$row17[0,4] = -0.7003938837051653
$row17[0,5] = 0.06973684323859711
$row17[0,6] = -0.1961638776409175
$ws.Range("B17:H17").Value = $row17
$ws.Range("I17").ClearContents()

$row18 = New-Object 'object[,]' 1,6
$row18[0,0] = -0.1460596688910241
$row18[0,1] = -0.08838812031281398
$row18[0,2] = 0.03557825052547153
$row18[0,3] = -0.8578398727093204
$row18[0,4] = -0.08770914576555788
$row18[0,5] = -0.3536098666450724
$ws.Range("B18:G18").Value = $row18
$ws.Range("H18").ClearContents()

$row19 = New-Object 'object[,]' 1,5
$row19[0,0] = 0.4122713362132568
$row19[0,1] = 0.5362377070515423
$row19[0,2] = -0.3571804161832495
$row19[0,3] = 0.4129503107605129
$row19[0,4] = 0.1470495898809984
$ws.Range("B19:F19").Value = $row19
$ws.Range("G19").ClearContents()

$row20 = New-Object 'object[,]' 1,4
$row20[0,0] = 0.2558798265160732
$row20[0,1] = -0.6375382967187186
$row20[0,2] = 0.1325924302250437
$row20[0,3] = -0.1333082906544708
$ws.Range("B20:E20").Value = $row20
$ws.Range("F20").ClearContents()

$row21 = New-Object 'object[,]' 1,3
$row21[0,0] = -0.465563438873707
$row21[0,1] = 0.3045672880700554
$row21[0,2] = 0.03866656719054083
$ws.Range("B21:D21").Value = $row21
$ws.Range("E21").ClearContents()

$row22 = New-Object 'object[,]' 1,2
$row22[0,0] = -0.001295251676276088
$row22[0,1] = -0.2671959725557906
$ws.Range("B22:C22").Value = $row22
$ws.Range("D22").ClearContents()

$row23 = New-Object 'object[,]' 1,1
$row23[0,0] = 0.3451339801314955
$ws.Range("B23:B23").Value = $row23
$ws.Range("C23").ClearContents()

$ws.Range("B24").ClearContents()
